# Update countries & provincias Spain
# Refresh COVID case data as of 13 de Abril de 2020 a las 17:22.
# Also resorts three countries (Letonia / Principado de Andorra / Republica
# de Chipre) whose case totals shuffled their rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 17:22"

# --- Updated totals for existing rows (same country, new numbers) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 562036
$ws.Cells.Item(4, 3).Value = 1736
$ws.Cells.Item(4, 4).Value = 33269
$ws.Cells.Item(4, 5).Value = 506604
$ws.Cells.Item(4, 6).Value = 11787
$ws.Cells.Item(4, 7).Value = 58
$ws.Cells.Item(4, 8).Value = 22163

# Row 8: Alemania
$ws.Cells.Item(8, 2).Value = 127916
$ws.Cells.Item(8, 3).Value = 62
$ws.Cells.Item(8, 5).Value = 60594

# Row 17: Brasil
$ws.Cells.Item(17, 2).Value = 22625
$ws.Cells.Item(17, 3).Value = 433
$ws.Cells.Item(17, 5).Value = 21207
$ws.Cells.Item(17, 7).Value = 22
$ws.Cells.Item(17, 8).Value = 1245

# Row 20: Austria
$ws.Cells.Item(20, 2).Value = 14018
$ws.Cells.Item(20, 3).Value = 73
$ws.Cells.Item(20, 5).Value = 6307

# Row 30: Polonia
$ws.Cells.Item(30, 2).Value = 6934
$ws.Cells.Item(30, 3).Value = 260
$ws.Cells.Item(30, 4).Value = 487
$ws.Cells.Item(30, 5).Value = 6202
$ws.Cells.Item(30, 7).Value = 13
$ws.Cells.Item(30, 8).Value = 245

# Row 36: Pakistan
$ws.Cells.Item(36, 2).Value = 5493
$ws.Cells.Item(36, 3).Value = 263
$ws.Cells.Item(36, 5).Value = 4305

# Row 54: Argentina
$ws.Cells.Item(54, 5).Value = 1643
$ws.Cells.Item(54, 7).Value = 7
$ws.Cells.Item(54, 8).Value = 97

# Row 56: Grecia
$ws.Cells.Item(56, 2).Value = 2145
$ws.Cells.Item(56, 3).Value = 31
$ws.Cells.Item(56, 5).Value = 1777
$ws.Cells.Item(56, 6).Value = 73
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 99

# Row 64: Irak
$ws.Cells.Item(64, 2).Value = 1378
$ws.Cells.Item(64, 3).Value = 26
$ws.Cells.Item(64, 4).Value = 717
$ws.Cells.Item(64, 5).Value = 583
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = 78

# Rows 87-89: Letonia / Principado de Andorra / Republica de Chipre swap
# order (Republica de Chipre's case count overtook the other two).
# Row 87 becomes Republica de Chipre (new data)
$ws.Cells.Item(87, 1).Value = "Republica de Chipre"
$ws.Cells.Item(87, 2).Value = 662
$ws.Cells.Item(87, 3).Value = 29
$ws.Cells.Item(87, 4).Value = 65
$ws.Cells.Item(87, 5).Value = 586
$ws.Cells.Item(87, 6).Value = 8
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 11

# Row 88 becomes Letonia (old row-87 data)
$ws.Cells.Item(88, 1).Value = "Letonia"
$ws.Cells.Item(88, 2).Value = 655
$ws.Cells.Item(88, 3).Value = 4
$ws.Cells.Item(88, 4).Value = 16
$ws.Cells.Item(88, 5).Value = 634
$ws.Cells.Item(88, 6).Value = 2
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 5

# Row 89 becomes Principado de Andorra (old row-88 data)
$ws.Cells.Item(89, 1).Value = "Principado de Andorra"
$ws.Cells.Item(89, 2).Value = 646
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 4).Value = 128
$ws.Cells.Item(89, 5).Value = 489
$ws.Cells.Item(89, 6).Value = 17
$ws.Cells.Item(89, 7).Value = 0
$ws.Cells.Item(89, 8).Value = 29

# Row 101: Reunion
$ws.Cells.Item(101, 2).Value = 391
$ws.Cells.Item(101, 3).Value = 2
$ws.Cells.Item(101, 5).Value = 351

# Row 119: Mayotte
$ws.Cells.Item(119, 2).Value = 207
$ws.Cells.Item(119, 3).Value = 11
$ws.Cells.Item(119, 5).Value = 145

# Row 132: Madagascar
$ws.Cells.Item(132, 4).Value = 21
$ws.Cells.Item(132, 5).Value = 85
